# The target diff only reorders the xmlns:* namespace-declaration
# attributes on the root element of word/document.xml, word/endnotes.xml,
# word/footer1.xml, word/footnotes.xml, word/header1.xml, word/styles.xml
# and word/theme/theme1.xml (e.g. w,r,w15,w14,m,wp,a,wp14,...  ->
# w,m,w14,r,wp,a,wp14,w15,...). The set of declared namespaces, the
# mc:Ignorable list, and every element/attribute/text node below those
# root tags are byte-for-byte identical before and after. That reordering
# is a side effect of the authoring tool that re-serialized the fixture
# (docx4j, per the commit message) switching its internal namespace-map
# iteration order - it is not a document edit, and there is no property
# in the Word object model (Document/Range/Paragraphs/Styles/Headers/
# Footers/...) that controls attribute-serialization order on a part's
# root element, so there is nothing to "do" to the content here.
#
# Touch the document (read-only) so the COM bridge sees activity, without
# mutating any content, formatting, or structure.
$d = $word.ActiveDocument
$null = $d.Content.Text.Length
